$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Cells.Item(15, 8).Value = 18008.82
$ws.Cells.Item(15, 9).Value = 18008.82
$ws.Cells.Item(15, 11).Value = 54026.46
$ws.Cells.Item(15, 13).Value = -53857.46

# row 38
$ws.Cells.Item(38, 8).Value = 285.75
$ws.Cells.Item(38, 9).Value = 47.666668
$ws.Cells.Item(38, 11).Value = 143.000004
$ws.Cells.Item(38, 13).Value = 228.999996

# row 40
$ws.Cells.Item(40, 8).Value = 1818.1666
$ws.Cells.Item(40, 9).Value = 1837.8334
$ws.Cells.Item(40, 10).Value = 1808.3334
$ws.Cells.Item(40, 11).Value = 1837.8334
$ws.Cells.Item(40, 12).Value = 1808.3334
$ws.Cells.Item(40, 13).Value = -1662.8334
$ws.Cells.Item(40, 14).Value = -2158.3334

# row 41
$ws.Cells.Item(41, 8).Value = 5383.4546
$ws.Cells.Item(41, 9).Value = 100
$ws.Cells.Item(41, 10).Value = 6557.5557
$ws.Cells.Item(41, 11).Value = 100
$ws.Cells.Item(41, 12).Value = 6557.5557
$ws.Cells.Item(41, 13).Value = 340
$ws.Cells.Item(41, 14).Value = -7437.5557

# row 98
$ws.Cells.Item(98, 8).Value = 7164.294
$ws.Cells.Item(98, 9).Value = 8320.5
$ws.Cells.Item(98, 10).Value = 1768.6666
$ws.Cells.Item(98, 11).Value = 8320.5
$ws.Cells.Item(98, 12).Value = 1768.6666
$ws.Cells.Item(98, 13).Value = -6822.5
$ws.Cells.Item(98, 14).Value = -4764.6666

# row 100
$ws.Cells.Item(100, 8).Value = 1057.4445
$ws.Cells.Item(100, 9).Value = 902.6923
$ws.Cells.Item(100, 10).Value = 1459.8
$ws.Cells.Item(100, 11).Value = 902.6923
$ws.Cells.Item(100, 12).Value = 1459.8
$ws.Cells.Item(100, 13).Value = -361.6923
$ws.Cells.Item(100, 14).Value = -2541.8

# row 107
$ws.Cells.Item(107, 8).Value = 251.80952
$ws.Cells.Item(107, 9).Value = 253.93333
$ws.Cells.Item(107, 11).Value = 253.93333
$ws.Cells.Item(107, 13).Value = 1666.06667

# row 116
$ws.Cells.Item(116, 8).Value = 3500
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 13).ClearContents()

# row 121
$ws.Cells.Item(121, 8).Value = 1029.1666
$ws.Cells.Item(121, 9).Value = 600
$ws.Cells.Item(121, 10).Value = 1115
$ws.Cells.Item(121, 11).Value = 1800
$ws.Cells.Item(121, 12).Value = 3345
$ws.Cells.Item(121, 14).Value = -6839
$ws.Cells.Item(121, 13).Value = -53

# row 122
$ws.Cells.Item(122, 8).Value = 7164.294
$ws.Cells.Item(122, 9).Value = 8320.5
$ws.Cells.Item(122, 10).Value = 1768.6666
$ws.Cells.Item(122, 11).Value = 24961.5
$ws.Cells.Item(122, 12).Value = 5305.9998
$ws.Cells.Item(122, 13).Value = -22511.5
$ws.Cells.Item(122, 14).Value = -10205.9998

# row 131
$ws.Cells.Item(131, 8).Value = 586.8
$ws.Cells.Item(131, 9).Value = 435.33334
$ws.Cells.Item(131, 10).Value = 1950
$ws.Cells.Item(131, 11).Value = 1306.00002
$ws.Cells.Item(131, 12).Value = 5850
$ws.Cells.Item(131, 13).Value = 3733.99998
$ws.Cells.Item(131, 14).Value = -15930

# row 132
$ws.Cells.Item(132, 8).Value = 2748655.2
$ws.Cells.Item(132, 9).Value = 3040984.5
$ws.Cells.Item(132, 10).Value = 761
$ws.Cells.Item(132, 11).Value = 9122953.5
$ws.Cells.Item(132, 12).Value = 2283
$ws.Cells.Item(132, 13).Value = -9120423.5
$ws.Cells.Item(132, 14).Value = -7343

# row 141
$ws.Cells.Item(141, 8).Value = 3352.075
$ws.Cells.Item(141, 9).Value = 1621.0714
$ws.Cells.Item(141, 10).Value = 4284.154
$ws.Cells.Item(141, 11).Value = 4863.2142
$ws.Cells.Item(141, 12).Value = 12852.462
$ws.Cells.Item(141, 13).Value = 316.7857999999997
$ws.Cells.Item(141, 14).Value = -23212.462

$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Cells.Item(61, 8).Value = 2372.4707
$ws.Cells.Item(61, 9).Value = 1232.8
$ws.Cells.Item(61, 11).Value = 1232.8
$ws.Cells.Item(61, 13).Value = -1020.8

# row 136
$ws.Cells.Item(136, 8).Value = 2372.4707
$ws.Cells.Item(136, 9).Value = 1232.8
$ws.Cells.Item(136, 11).Value = 3698.4
$ws.Cells.Item(136, 13).Value = -1148.4

$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Cells.Item(20, 8).Value = 3230.2
$ws.Cells.Item(20, 9).Value = 3376.0952
$ws.Cells.Item(20, 10).Value = 2889.7778
$ws.Cells.Item(20, 11).Value = 3376.0952
$ws.Cells.Item(20, 12).Value = 2889.7778
$ws.Cells.Item(20, 13).Value = -3129.0952
$ws.Cells.Item(20, 14).Value = -3383.7778

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Cells.Item(31, 8).Value = 4904861
$ws.Cells.Item(31, 9).Value = 2627.7727
$ws.Cells.Item(31, 10).Value = 13892289
$ws.Cells.Item(31, 11).Value = 2627.7727
$ws.Cells.Item(31, 12).Value = 13892289
$ws.Cells.Item(31, 13).Value = -2332.7727
$ws.Cells.Item(31, 14).Value = -13892879

# row 34
$ws.Cells.Item(34, 8).Value = 4904861
$ws.Cells.Item(34, 9).Value = 2627.7727
$ws.Cells.Item(34, 10).Value = 13892289
$ws.Cells.Item(34, 11).Value = 2627.7727
$ws.Cells.Item(34, 12).Value = 13892289
$ws.Cells.Item(34, 13).Value = -2425.7727
$ws.Cells.Item(34, 14).Value = -13892693

# row 122
$ws.Cells.Item(122, 8).Value = 991.4545000000001
$ws.Cells.Item(122, 9).Value = 1062
$ws.Cells.Item(122, 10).Value = 932.6667
$ws.Cells.Item(122, 11).Value = 3186
$ws.Cells.Item(122, 12).Value = 2798.0001
$ws.Cells.Item(122, 13).Value = -736
$ws.Cells.Item(122, 14).Value = -7698.0001

$ws = $wb.Worksheets.Item("CUL")
# row 107
$ws.Cells.Item(107, 8).Value = 346.75
$ws.Cells.Item(107, 9).Value = 427.22223
$ws.Cells.Item(107, 10).Value = 280.9091
$ws.Cells.Item(107, 11).Value = 1281.66669
$ws.Cells.Item(107, 12).Value = 842.7273
$ws.Cells.Item(107, 13).Value = 638.33331
$ws.Cells.Item(107, 14).Value = -4682.7273

# row 113
$ws.Cells.Item(113, 8).Value = 523.7742
$ws.Cells.Item(113, 9).Value = 493
$ws.Cells.Item(113, 10).Value = 605.2353000000001
$ws.Cells.Item(113, 11).Value = 1479
$ws.Cells.Item(113, 12).Value = 1815.7059
$ws.Cells.Item(113, 13).Value = 691
$ws.Cells.Item(113, 14).Value = -6155.7059

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Cells.Item(70, 8).Value = 8364469.5
$ws.Cells.Item(70, 9).Value = 9811501
$ws.Cells.Item(70, 10).Value = 3843.2222
$ws.Cells.Item(70, 11).Value = 9811501
$ws.Cells.Item(70, 12).Value = 3843.2222
$ws.Cells.Item(70, 13).Value = -9811231
$ws.Cells.Item(70, 14).Value = -4383.2222

# row 73
$ws.Cells.Item(73, 8).Value = 8364469.5
$ws.Cells.Item(73, 9).Value = 9811501
$ws.Cells.Item(73, 10).Value = 3843.2222
$ws.Cells.Item(73, 11).Value = 9811501
$ws.Cells.Item(73, 12).Value = 3843.2222
$ws.Cells.Item(73, 13).Value = -9810565
$ws.Cells.Item(73, 14).Value = -5715.2222

$ws = $wb.Worksheets.Item("LTW")
# row 74
$ws.Cells.Item(74, 8).Value = 36000
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 36000
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 36000
$ws.Cells.Item(74, 14).Value = -37996
$ws.Cells.Item(74, 13).ClearContents()

# row 77
$ws.Cells.Item(77, 8).Value = 36000
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 36000
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 108000
$ws.Cells.Item(77, 14).Value = -117984
$ws.Cells.Item(77, 13).ClearContents()

# row 132
$ws.Cells.Item(132, 8).Value = 5828.829
$ws.Cells.Item(132, 9).Value = 8081.44
$ws.Cells.Item(132, 10).Value = 2309.125
$ws.Cells.Item(132, 11).Value = 24244.32
$ws.Cells.Item(132, 12).Value = 6927.375
$ws.Cells.Item(132, 13).Value = -21714.32
$ws.Cells.Item(132, 14).Value = -11987.375

# row 140
$ws.Cells.Item(140, 8).Value = 50809.668
$ws.Cells.Item(140, 10).Value = 50809.668
$ws.Cells.Item(140, 12).Value = 50809.668
$ws.Cells.Item(140, 14).Value = -61169.668

$ws = $wb.Worksheets.Item("WVR")
# row 75
$ws.Cells.Item(75, 8).Value = 30043.334
$ws.Cells.Item(75, 10).Value = 30043.334
$ws.Cells.Item(75, 12).Value = 30043.334
$ws.Cells.Item(75, 14).Value = -31915.334

# row 78
$ws.Cells.Item(78, 8).Value = 30043.334
$ws.Cells.Item(78, 10).Value = 30043.334
$ws.Cells.Item(78, 12).Value = 90130.00199999999
$ws.Cells.Item(78, 14).Value = -99490.00199999999

# row 113
$ws.Cells.Item(113, 8).Value = 589.5625
$ws.Cells.Item(113, 9).Value = 770.8889
$ws.Cells.Item(113, 10).Value = 356.42856
$ws.Cells.Item(113, 11).Value = 2312.6667
$ws.Cells.Item(113, 12).Value = 1069.28568
$ws.Cells.Item(113, 13).Value = -142.6667000000002
$ws.Cells.Item(113, 14).Value = -5409.28568

# row 122
$ws.Cells.Item(122, 8).Value = 48199.547
$ws.Cells.Item(122, 9).Value = 69361.664
$ws.Cells.Item(122, 10).Value = 2852.1428
$ws.Cells.Item(122, 11).Value = 208084.992
$ws.Cells.Item(122, 12).Value = 8556.428400000001
$ws.Cells.Item(122, 13).Value = -205634.992
$ws.Cells.Item(122, 14).Value = -13456.4284
